$d = $word.ActiveDocument
$anchor = $d.Content
$anchor.Find.Execute("Sales Performance:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchRng = $d.Range($anchor.End, $anchor.End + 10)
$searchRng.Find.Execute(" done", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchRng.Collapse(0)
$searchRng.InsertAfter(", can be shown on Looker Studio")
Write-Output $searchRng.Text
$searchRng.Font.NameAscii = "Times New Roman"
$searchRng.Font.Name = "Times New Roman"
$searchRng.Font.NameBi = "Times New Roman"
$searchRng.LanguageID = "en-US"
